$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 55, pushing existing rows 55-91 down to 56-92.
# The inserted row inherits formatting from the row above it (row 54) by
# default in Excel, so explicitly copy row 55's (pre-shift) formatting by
# duplicating its original values first, then adjust the two changed cells.
$ws.Rows("55:55").Insert(-4121)  # xlShiftDown

# New row 55 is a copy of what used to be row 55 (now row 56), with an
# updated date (D) and volume (M).
$ws.Range("A55").Value = 1
$ws.Range("B55").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C55").Value = "Arica y Parinacota"
$ws.Range("D55").Value = 44762
$ws.Range("E55").Value = 15
$ws.Range("F55").Value = "Fruta"
$ws.Range("G55").Value = 100102
$ws.Range("H55").Value = "Cítricos"
$ws.Range("I55").Value = 100102005
$ws.Range("J55").Value = "Naranja"
$ws.Range("K55").Value = "Fukumoto"
$ws.Range("L55").Value = "Segunda"
$ws.Range("M55").Value = 300
$ws.Range("N55").Value = 600
$ws.Range("O55").Value = 650
$ws.Range("P55").Value = 625
$ws.Range("Q55").Value = "$/kilo (en caja de 20 kilos)"
$ws.Range("R55").Value = "Región de Coquimbo"
$ws.Range("S55").Value = 625
$ws.Range("T55").Value = 1
